# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column O holds the "municipio-nombre" field metadata block (rows 2-4).
# It is re-curated to be treated the same way as the existing
# "provincia-nombre" (P) and "comarca-nombre" (Q) reference-area
# dimensions: a sdmx-dimension:refArea / dim / URI-Municipio triple,
# instead of the old iaest-measure:municipio-nombre / medida / xsd:int
# triple.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("O2").Value = "sdmx-dimension:refArea"
$ws.Range("O3").Value = "dim"
$ws.Range("O4").Value = "URI-Municipio"
